$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.984.41'
$ws.Range('E2').Value = '  +3.19%  '
$ws.Range('D3').Value = '3.472.90'
$ws.Range('E3').Value = '  +3.53%  '
$ws.Range('E4').Value = '  +0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '408.09'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.37%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '132.07'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +17.22%  '
$ws.Range('D7').Value = '3.465.06'
$ws.Range('E7').Value = '  +3.48%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.605'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +3.09%  '
$ws.Range('E9').Value = '  +0.02%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.694'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +8.54%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.134'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +34.28%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '43.63'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +9.54%  '
$ws.Range('E13').Value = '  -0.84%  '
$ws.Range('D14').Value = '4.021.05'
$ws.Range('E14').Value = '  +3.59%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '8.85'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +5.05%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '20.14'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('D17').Value = '3.478.57'
$ws.Range('E17').Value = '  +4.03%  '
$ws.Range('D18').Value = '62.995.71'
$ws.Range('E18').Value = '  +3.69%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('E20').Value = '  +1.37%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '0.0000144'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +29.03%  '
$ws.Range('E22').Value = '  -1.30%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '82.41'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +9.81%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '13.15'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.24%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '312.54'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.95%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '3.18'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  +5.95%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '8.21'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +3.16%  '
$ws.Range('E29').Value = '  -1.12%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.37'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -2.60%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.57'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('E32').Value = '  +3.01%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '44.19'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +13.00%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '11.83'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +3.42%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.59'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -1.73%  '
$ws.Range('E36').Value = '  -0.07%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.0495'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -2.91%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '52.67'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('E39').Value = '  +4.50%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  -3.04%  '
$ws.Range('E42').Value = '  +2.67%  '
$ws.Range('E43').Value = '  +4.17%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '137.35'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.45%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '17.56'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +4.00%  '
$ws.Range('E46').Value = '  -2.85%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '4.00'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +0.85%  '
$ws.Range('E48').Value = '  -0.35%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '22.19'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('D50').Value = '3.819.21'
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('D51').Value = '2.185.54'
$ws.Range('E51').Value = '  +0.46%  '
